$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# --- Crime statistics table updates (rows 15-33) ---
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 133.333333333333
$ws.Range("L15").Value = 250
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = -65
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = 4.081632653061
$ws.Range("L16").Value = 10.869565217391
$ws.Range("M16").Value = -51.428571428571
$ws.Range("N16").Value = -90.810810810810
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -30.434782608695
$ws.Range("I17").Value = 114
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = 15.151515151515
$ws.Range("L17").Value = 26.666666666666
$ws.Range("M17").Value = 80.952380952380
$ws.Range("N17").Value = -51.898734177215
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = -29.032258064516
$ws.Range("L18").Value = 2.325581395348
$ws.Range("M18").Value = -38.888888888888
$ws.Range("N18").Value = -90.350877192982
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -29.032258064516
$ws.Range("I19").Value = 165
$ws.Range("J19").Value = 141
$ws.Range("K19").Value = 17.021276595744
$ws.Range("L19").Value = 4.430379746835
$ws.Range("M19").Value = -24.657534246575
$ws.Range("N19").Value = -36.781609195402
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 11.111111111111
$ws.Range("I20").Value = 53
$ws.Range("J20").Value = 67
$ws.Range("K20").Value = -20.895522388059
$ws.Range("L20").Value = 60.606060606060
$ws.Range("M20").Value = 1.923076923076
$ws.Range("N20").Value = -85.112359550561
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = -21.951219512195
$ws.Range("I21").Value = 435
$ws.Range("J21").Value = 422
$ws.Range("K21").Value = 3.080568720379
$ws.Range("L21").Value = 16.310160427807
$ws.Range("M21").Value = -16.184971098265
$ws.Range("N21").Value = -77.117306680694
$ws.Range("M22").Value = -75
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 11.111111111111
$ws.Range("I23").Value = 69
$ws.Range("J23").Value = 54
$ws.Range("K23").Value = 27.777777777777
$ws.Range("L23").Value = 25.454545454545
$ws.Range("M23").Value = 25.454545454545
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 55
$ws.Range("H24").Value = 3.636363636363
$ws.Range("I24").Value = 333
$ws.Range("J24").Value = 304
$ws.Range("K24").Value = 9.539473684210
$ws.Range("L24").Value = -8.767123287671
$ws.Range("M24").Value = -38.218923933209
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 88.888888888888
$ws.Range("I25").Value = 86
$ws.Range("J25").Value = 76
$ws.Range("K25").Value = 13.157894736842
$ws.Range("L25").Value = -9.473684210526
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = 60.869565217391
$ws.Range("I26").Value = 203
$ws.Range("J26").Value = 187
$ws.Range("K26").Value = 8.556149732620
$ws.Range("L26").Value = 16.666666666666
$ws.Range("M26").Value = 22.289156626506
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 133.333333333333
$ws.Range("L27").Value = 75
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 17
$ws.Range("K28").Value = 58.823529411764
$ws.Range("L28").Value = 17.391304347826
$ws.Range("N29").Value = -90.163934426229
$ws.Range("N30").Value = -88.461538461538
$ws.Range("I33").Value = 3

# --- Cells whose type changes (text "N/A" <-> number) need a format fix ---
# so the cell reuses the existing number-style / text-style slot exactly
# like Excel would (copy format only, value already set).
$ws.Range("D15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = 100
$ws.Range("K14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").Value = 100
$ws.Range("K14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("G28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H28").Value = 200
$ws.Range("K14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("D31").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("G31").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("H31").Value = -100
$ws.Range("K14").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("J31").Value = 1
$ws.Range("I14").Copy()
$ws.Range("J31").PasteSpecial(-4122)
$ws.Range("K31").Value = 400
$ws.Range("K14").Copy()
$ws.Range("K31").PasteSpecial(-4122)
$ws.Range("C33").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("F33").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F33").PasteSpecial(-4122)

$excel.CutCopyMode = 0
